$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.836.57'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.50%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.806.82'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.03%  '

# Row 4
$ws.Range('E4').Value = '  -0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '353.01'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.38%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '112.31'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.71%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.558'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.36%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('E9').Value = '  +7.45%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.28'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.35%  '

# Row 11
$ws.Range('E11').Value = '  -0.30%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0840'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.28%  '

# Row 13
$ws.Range('E13').Value = '  +1.41%  '

# Row 14
$ws.Range('E14').Value = '  +4.55%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.245.77'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.25%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.804.94'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.38%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.950'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.33%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.830.07'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.68%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.63'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.27%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.28'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +8.34%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.54'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +4.84%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0974'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.34%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.35'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.93%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '267.42'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.77%  '

# Row 25
$ws.Range('E25').Value = '  +2.19%  '

# Row 26
$ws.Range('E26').Value = '  +0.11%  '

# Row 27
$ws.Range('E27').Value = '  +1.57%  '

# Row 28
$ws.Range('E28').Value = '  +0.44%  '

# Row 29
$ws.Range('E29').Value = '  +13.77%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '10.40'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +3.93%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.29'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.35%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '52.19'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.14%  '

# Row 33
$ws.Range('E33').Value = '  +2.44%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0896'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +8.78%  '

# Row 35
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0451'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.02%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.52'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +6.83%  '

# Row 37
$ws.Range('E37').Value = '  -0.02%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '19.02'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.57%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.17'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.57%  '

# Row 40
$ws.Range('E40').Value = '  +3.88%  '

# Row 41
$ws.Range('E41').Value = '  +2.45%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.52'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.79%  '

# Row 43
$ws.Range('E43').Value = '  +1.60%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '119.81'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.93%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.92'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.48%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.53'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +10.49%  '

# Row 47
$ws.Range('E47').Value = '  +9.23%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.113.90'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.66%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.988'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +8.22%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.51'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.55%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.38'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +8.67%  '
